$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 539 (shifts existing rows 539..667 down to 540..668)
$ws.Rows(539).Insert()

# Populate the newly inserted row 539 with the new record
$ws.Cells.Item(539, 1).Value = 5
$ws.Cells.Item(539, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(539, 3).Value = "Maule"
$ws.Cells.Item(539, 4).Value = [DateTime]"2023-10-05"
$ws.Cells.Item(539, 5).Value = 7
$ws.Cells.Item(539, 6).Value = 100114014
$ws.Cells.Item(539, 7).Value = "Betarraga"
$ws.Cells.Item(539, 8).Value = "Sin especificar"
$ws.Cells.Item(539, 9).Value = "Primera"
$ws.Cells.Item(539, 10).Value = 4000
$ws.Cells.Item(539, 11).Value = 500
$ws.Cells.Item(539, 12).Value = 500
$ws.Cells.Item(539, 13).Value = 500
$ws.Cells.Item(539, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(539, 15).Value = "Región del Maule"
$ws.Cells.Item(539, 16).Value = 100
$ws.Cells.Item(539, 17).Value = 5
$ws.Cells.Item(539, 18).Value = "Hortaliza"
